$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C12").Value = 0.5816
$ws.Range("C13").Value = 0.7846
$ws.Range("C14").Value = 0.4931
$ws.Range("C15").Value = 0.5071
$ws.Range("C16").Value = 0.8415
$ws.Range("C17").Value = 0.7586000000000001
$ws.Range("C18").Value = 0.6773
$ws.Range("C19").Value = 0.7846
$ws.Range("C20").Value = 0.7586000000000001
$ws.Range("C21").Value = 0.922
$ws.Range("C22").Value = 0.3943
$ws.Range("C23").Value = 0.9552
$ws.Range("C24").Value = 0.8369
$ws.Range("C25").Value = 0.9228
$ws.Range("C26").Value = 1
$ws.Range("C27").Value = 0.7695
$ws.Range("C28").Value = 0.7967
$ws.Range("C29").Value = 0.6414
$ws.Range("C30").Value = 0.6844
$ws.Range("C31").Value = 0.626
$ws.Range("C32").Value = 0.869
$ws.Range("C33").Value = 0.4965
$ws.Range("C34").Value = 0.7724
$ws.Range("C35").Value = 0.9517
$ws.Range("C36").Value = 0.3723
$ws.Range("C37").Value = 0.4675
$ws.Range("C38").Value = 0.6516999999999999
$ws.Range("C39").Value = 0.7589
$ws.Range("C40").Value = 0.4797
$ws.Range("C41").Value = 0.7171999999999999
$ws.Range("C42").Value = 0.9043
$ws.Range("C43").Value = 0.4024
$ws.Range("C44").Value = 0.5793
$ws.Range("C45").Value = 0.8723
$ws.Range("C46").Value = 0.7886
